$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to retain exact text formatting (avoid numeric coercion
# of values like "1.000", "0.7370", "30.397.23") by marking the range as
# Text before writing, then resetting the style back to Normal afterwards
# so no stray number-format style is left behind.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = "30.397.23"
$ws.Range("E2").Value = "  -2.31%  "
$ws.Range("D3").Value = "1.915.60"
$ws.Range("E3").Value = "  +0.95%  "
$ws.Range("D4").Value = "0.9999"
$ws.Range("E4").Value = "  +0.39%  "
$ws.Range("D5").Value = "241.12"
$ws.Range("E5").Value = "  +0.83%  "
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.47%  "
$ws.Range("D7").Value = "0.4693"
$ws.Range("E7").Value = "  -2.26%  "
$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").Value = "0.2844"
$ws.Range("E8").Value = "  -0.77%  "
$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").Value = "0.06951"
$ws.Range("E9").Value = "  +5.97%  "
$ws.Range("B10").Value = "Litecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D10").Value = "106.29"
$ws.Range("E10").Value = "  +10.65%  "
$ws.Range("B11").Value = "Solana"
$ws.Range("C11").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D11").Value = "18.07"
$ws.Range("E11").Value = "  -4.21%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.901.08"
$ws.Range("E12").Value = "  +0.06%  "
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").Value = "0.07637"
$ws.Range("E13").Value = "  +1.45%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "5.175"
$ws.Range("E14").Value = "  +0.36%  "
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").Value = "0.6555"
$ws.Range("E15").Value = "  -0.19%  "
$ws.Range("B16").Value = "BitcoinCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D16").Value = "285.98"
$ws.Range("E16").Value = "  -4.32%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "30.409.38"
$ws.Range("E17").Value = "  -2.12%  "
$ws.Range("B18").Value = "Avalanche"
$ws.Range("C18").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D18").Value = "12.99"
$ws.Range("E18").Value = "  -1.51%  "
$ws.Range("D19").Value = "1.001"
$ws.Range("E19").Value = "  +0.41%  "
$ws.Range("D20").Value = "0.000007596"
$ws.Range("E20").Value = "  +0.13%  "
$ws.Range("B21").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C21").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D21").Value = "2.146.42"
$ws.Range("E21").Value = "  +1.57%  "
$ws.Range("B22").Value = "BinanceUSD"
$ws.Range("C22").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  +0.44%  "
$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").Value = "5.223"
$ws.Range("E23").Value = "  +0.12%  "
$ws.Range("B24").Value = "Chainlink"
$ws.Range("C24").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D24").Value = "6.184"
$ws.Range("E24").Value = "  +0.22%  "
$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").Value = "168.13"
$ws.Range("E25").Value = "  -0.17%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").Value = "9.236"
$ws.Range("E26").Value = "  -1.12%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "21.05"
$ws.Range("E27").Value = "  +7.00%  "
$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").Value = "2.032"
$ws.Range("E28").Value = "  +2.74%  "
$ws.Range("B29").Value = "Stellar"
$ws.Range("C29").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D29").Value = "0.1071"
$ws.Range("E29").Value = "  +1.28%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "1.368"
$ws.Range("E30").Value = "  -0.36%  "
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").Value = "4.128"
$ws.Range("E31").Value = "  -0.99%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "3.952"
$ws.Range("E32").Value = "  -1.23%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "0.05054"
$ws.Range("E33").Value = "  +0.65%  "
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").Value = "0.7370"
$ws.Range("E34").Value = "  +1.06%  "
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value = "1.144"
$ws.Range("E35").Value = "  -3.91%  "
$ws.Range("B36").Value = "Frax"
$ws.Range("C36").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D36").Value = "0.9996"
$ws.Range("E36").Value = "  +0.29%  "
$ws.Range("D37").Value = "2.713"
$ws.Range("E37").Value = "  +0.41%  "
$ws.Range("E38").Value = "  +2.93%  "
$ws.Range("D39").Value = "2.667"
$ws.Range("E39").Value = "  -2.10%  "
$ws.Range("D40").Value = "2.049"
$ws.Range("E40").Value = "  -1.18%  "
$ws.Range("D41").Value = "108.44"
$ws.Range("E41").Value = "  +0.76%  "
$ws.Range("D42").Value = "0.8729"
$ws.Range("E42").Value = "  -3.35%  "
$ws.Range("D43").Value = "5.831"
$ws.Range("E43").Value = "  +4.07%  "
$ws.Range("D44").Value = "1.0000"
$ws.Range("E44").Value = "  +0.39%  "
$ws.Range("B45").Value = "TheSandbox"
$ws.Range("C45").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D45").Value = "0.4196"
$ws.Range("E45").Value = "  -1.25%  "
$ws.Range("B46").Value = "BitcoinSV"
$ws.Range("C46").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D46").Value = "52.23"
$ws.Range("E46").Value = "  +23.48%  "
$ws.Range("D47").Value = "67.11"
$ws.Range("E47").Value = "  +1.65%  "
$ws.Range("D48").Value = "7.126"
$ws.Range("E48").Value = "  -4.32%  "
$ws.Range("D49").Value = "9.164"
$ws.Range("E49").Value = "  +3.01%  "
$ws.Range("D50").Value = "0.1202"
$ws.Range("E50").Value = "  -2.85%  "
$ws.Range("D51").Value = "34.57"
$ws.Range("E51").Value = "  -0.55%  "

$dRange.Style = "Normal"

